$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Summary"
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = 836.76
$wsSummary.Range("E2").Value = 9163.24
$wsSummary.Range("F2").Value = 849.4
# Materialize an (empty, unstyled) G2 cell so the row/column extent grows to G,
# matching the new column G that appears in the edited workbook.
$wsSummary.Range("G2").Style = "Normal"

$wsSummary.Range("A3").Value = 561.21
$wsSummary.Range("E3").Value = 510.25
$wsSummary.Range("F3").Value = 38.32

$wsSummary.Range("A5").Value = 16.86
$wsSummary.Range("B5").Value = 8.3699999999999992

# Widen column B (used to be best-fit, now a fixed custom width)
$wsSummary.Columns.Item(2).ColumnWidth = 10.022135416666666

# ---------------------------------------------------------------------
# Sheet "Repayment schedule"
# ---------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

$wsRepay.Range("G3").Value = 4163.24

$wsRepay.Range("F5").Value = 849.4
$wsRepay.Range("G5").Value = 8313.84
$wsRepay.Range("H5").Value = 38.32

$wsRepay.Range("F6").Value = 794.33
$wsRepay.Range("G6").Value = 7519.51
$wsRepay.Range("H6").Value = 93.39

$wsRepay.Range("F7").Value = 813.55
$wsRepay.Range("G7").Value = 6705.96
$wsRepay.Range("H7").Value = 74.17

$wsRepay.Range("F8").Value = 819.37
$wsRepay.Range("G8").Value = 5886.59
$wsRepay.Range("H8").Value = 68.349999999999994

$wsRepay.Range("F9").Value = 829.66
$wsRepay.Range("G9").Value = 5056.93
$wsRepay.Range("H9").Value = 58.06

$wsRepay.Range("F10").Value = 836.18
$wsRepay.Range("G10").Value = 4220.75
$wsRepay.Range("H10").Value = 51.54

$wsRepay.Range("F11").Value = 844.7
$wsRepay.Range("G11").Value = 3376.05
$wsRepay.Range("H11").Value = 43.02

$wsRepay.Range("F12").Value = 854.42
$wsRepay.Range("G12").Value = 2521.63
$wsRepay.Range("H12").Value = 33.299999999999997

$wsRepay.Range("F13").Value = 862.02
$wsRepay.Range("G13").Value = 1659.61
$wsRepay.Range("H13").Value = 25.7

$wsRepay.Range("F14").Value = 871.35
$wsRepay.Range("G14").Value = 788.26
$wsRepay.Range("H14").Value = 16.37

$wsRepay.Range("F15").Value = 788.26
$wsRepay.Range("H15").Value = 8.0299999999999994
$wsRepay.Range("K15").Value = 796.29
$wsRepay.Range("P15").Value = 796.29

# Widen column L (used to be best-fit, now a fixed custom width)
$wsRepay.Columns.Item(12).ColumnWidth = 6.307291666666667

# ---------------------------------------------------------------------
# Sheet "Transactions"
# ---------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")

$wsTrans.Range("A2").Value = 881
$wsTrans.Range("J2").NumberFormat = "#,##0.00"
$wsTrans.Range("J2").Value = 9163.24

$wsTrans.Range("A3").Value = 875
$wsTrans.Range("J3").NumberFormat = "#,##0.00"
$wsTrans.Range("J3").Value = 4163.24

$wsTrans.Range("A4").Value = 858

# Widen column A slightly
$wsTrans.Columns.Item(1).ColumnWidth = 3.1666666666666665

# ---------------------------------------------------------------------
# Selections (match final cursor position captured in each sheet)
# ---------------------------------------------------------------------
$wsSummary.Select()
$wsSummary.Range("D5").Select()

$wsTrans.Select()
$wsTrans.Range("D4").Select()

# "Repayment schedule" ends up the active/visible tab, as in the source file
$wsRepay.Select()
$wsRepay.Range("F13").Select()
